$d = $word.ActiveDocument

# Locate the Subtitle paragraph ("Identity, desistance and the experience
# of imprisonment") so a new "Author" paragraph ("Ben Jarman") can be
# inserted immediately after it, before the Date paragraph.
$targetText = "Identity, desistance and the experience of imprisonment"

$found = $false
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    $t = $t.TrimEnd([char]13, [char]10, [char]7)
    if ($t -eq $targetText) {
        # Create a new empty paragraph right after this one.
        $p.Range.InsertParagraphAfter()

        # The newly created paragraph is now the next paragraph after $p.
        $newPara = $p.Next()

        $authorXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
            '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
            '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
            '<pkg:xmlData>' +
            '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
            '<w:body>' +
            '<w:p>' +
            '<w:pPr><w:pStyle w:val="Author"/></w:pPr>' +
            '<w:r><w:t xml:space="preserve">Ben Jarman</w:t></w:r>' +
            '</w:p>' +
            '</w:body>' +
            '</w:document>' +
            '</pkg:xmlData></pkg:part></pkg:package>'

        $newPara.Range.InsertXML($authorXml)

        $found = $true
        break
    }
}

if (-not $found) {
    throw "Could not locate the subtitle paragraph to insert the author after."
}
